$wb = $excel.ActiveWorkbook

# The new "Turkey" sheet is a trimmed-down copy of the existing "Spain"
# sheet (same layout/template, fewer repeater rows, different market name
# and ticket reference). Start by duplicating Spain, then adjust.
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)

$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Spain's template listed every repeater model (MXR..P32DR); Turkey's
# sheet only keeps PR1DS/PR8AS (plus the trailing Wg/Repeaters rows), so
# remove the rows for the models that don't apply - old rows 8:17.
$turkey.Rows("8:17").Delete()

# Update the market name + Jira/NGC ticket reference for Turkey.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3309/"

# Match the column widths used on the Turkey tab.
$turkey.Columns("B").ColumnWidth = 15.1666666666667
$turkey.Columns("D").ColumnWidth = 20.1666666666667

# Rows 3-5 carried Spain's taller (wrapped) row height; Turkey's sheet
# uses the default row height instead, so drop the override.
$turkey.Rows("3:5").AutoFit()

# Restore Spain's own selection (it's no longer the active tab) and set
# Turkey as the newly active/selected tab with its own selection.
[void]$spain.Activate()
[void]$spain.Range("A20:A21").Select()

[void]$turkey.Activate()
[void]$turkey.Range("G15").Select()
